$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2025-06-17 Tuesday" "2025-06-18 Wednesday"

Replace-Text "371÷2=185, 1" "600÷5=120, 0"
Replace-Text "734÷6=122, 2" "755÷8=94, 3"
Replace-Text "120÷6=20, 0" "692÷4=173, 0"
Replace-Text "353÷2=176, 1" "883÷8=110, 3"
Replace-Text "853÷7=121, 6" "957÷3=319, 0"

Replace-Text "334÷8=41, 6" "157÷7=22, 3"
Replace-Text "816÷8=102, 0" "829÷3=276, 1"
Replace-Text "312÷2=156, 0" "685÷9=76, 1"
Replace-Text "626÷2=313, 0" "601÷3=200, 1"
Replace-Text "803÷6=133, 5" "737÷9=81, 8"

Replace-Text "129÷7=18, 3" "940÷3=313, 1"
Replace-Text "264÷2=132, 0" "323÷2=161, 1"
Replace-Text "607÷8=75, 7" "905÷2=452, 1"
Replace-Text "636÷9=70, 6" "599÷6=99, 5"
Replace-Text "524÷6=87, 2" "580÷3=193, 1"

Replace-Text "352÷4=88, 0" "753÷5=150, 3"
Replace-Text "392÷5=78, 2" "474÷8=59, 2"
Replace-Text "992÷4=248, 0" "533÷6=88, 5"
Replace-Text "987÷3=329, 0" "324÷9=36, 0"
Replace-Text "674÷3=224, 2" "172÷2=86, 0"

Replace-Text "970÷6=161, 4" "776÷7=110, 6"
Replace-Text "231÷3=77, 0" "433÷6=72, 1"
Replace-Text "670÷6=111, 4" "606÷9=67, 3"
Replace-Text "640÷4=160, 0" "197÷9=21, 8"
Replace-Text "933÷6=155, 3" "453÷9=50, 3"

Write-Output "Done"
